$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 68 (pushes existing rows 68-101 down to 70-103)
$ws.Rows.Item(68).Resize(2).Insert()

# New row 68: Sin especificar / Especial / 400 / 5500 / 6000 / 5750 / $/bandeja 4 kilos / Perú / 1438 / 4
$ws.Cells.Item(68, 1).Value = 1
$ws.Cells.Item(68, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(68, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(68, 4).Value = 44523
$ws.Cells.Item(68, 5).Value = 15
$ws.Cells.Item(68, 6).Value = "Fruta"
$ws.Cells.Item(68, 7).Value = 100108
$ws.Cells.Item(68, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(68, 9).Value = 100108002
$ws.Cells.Item(68, 10).Value = "Mango"
$ws.Cells.Item(68, 11).Value = "Sin especificar"
$ws.Cells.Item(68, 12).Value = "Especial"
$ws.Cells.Item(68, 13).Value = 400
$ws.Cells.Item(68, 14).Value = 5500
$ws.Cells.Item(68, 15).Value = 6000
$ws.Cells.Item(68, 16).Value = 5750
$ws.Cells.Item(68, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(68, 18).Value = "Perú"
$ws.Cells.Item(68, 19).Value = 1438
$ws.Cells.Item(68, 20).Value = 4

# New row 69: Sin especificar / Primera / 400 / 5500 / 6000 / 5750 / $/bandeja 4 kilos / Perú / 1438 / 4
$ws.Cells.Item(69, 1).Value = 1
$ws.Cells.Item(69, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(69, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(69, 4).Value = 44523
$ws.Cells.Item(69, 5).Value = 15
$ws.Cells.Item(69, 6).Value = "Fruta"
$ws.Cells.Item(69, 7).Value = 100108
$ws.Cells.Item(69, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(69, 9).Value = 100108002
$ws.Cells.Item(69, 10).Value = "Mango"
$ws.Cells.Item(69, 11).Value = "Sin especificar"
$ws.Cells.Item(69, 12).Value = "Primera"
$ws.Cells.Item(69, 13).Value = 400
$ws.Cells.Item(69, 14).Value = 5500
$ws.Cells.Item(69, 15).Value = 6000
$ws.Cells.Item(69, 16).Value = 5750
$ws.Cells.Item(69, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(69, 18).Value = "Perú"
$ws.Cells.Item(69, 19).Value = 1438
$ws.Cells.Item(69, 20).Value = 4

# Match date-format style used by column D on the other rows
$ws.Cells.Item(68, 4).NumberFormat = $ws.Cells.Item(70, 4).NumberFormat
$ws.Cells.Item(69, 4).NumberFormat = $ws.Cells.Item(70, 4).NumberFormat
